$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source cells are plain text (e.g. "321.16", "7.37%") rather than
# numbers/percentages. Temporarily mark each target cell as Text before
# assigning so Excel does not auto-convert the value, then restore the
# cell style so formatting matches the original (unstyled) cells.
function Set-TextValue($ref, $value) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue "D2" "321.16"
Set-TextValue "E2" "7.37%"
Set-TextValue "D3" "48.84"
Set-TextValue "E3" "15.75%"
Set-TextValue "D4" "5.263"
Set-TextValue "E4" "4.93%"
Set-TextValue "D5" "0.08110"
Set-TextValue "E5" "7.47%"
Set-TextValue "D6" "4.590"
Set-TextValue "E6" "5.25%"
Set-TextValue "D7" "1.643"
Set-TextValue "E7" "2.45%"
Set-TextValue "D8" "1.206"
Set-TextValue "E8" "28.80%"
Set-TextValue "D9" "0.1292"
Set-TextValue "E9" "8.76%"
Set-TextValue "D10" "0.1950"
Set-TextValue "E10" "6.00%"
Set-TextValue "D11" "0.09500"
Set-TextValue "E11" "4.79%"
Set-TextValue "D12" "0.04624"
Set-TextValue "E12" "11.13%"
Set-TextValue "D13" "0.1051"
Set-TextValue "E13" "0.34%"
Set-TextValue "D14" "0.001334"
Set-TextValue "E14" "3.08%"
Set-TextValue "D15" "0.04161"
Set-TextValue "E15" "1.92%"
Set-TextValue "D16" "0.005947"
Set-TextValue "E16" "2.60%"
Set-TextValue "D17" "3.342"
Set-TextValue "E17" "0.05%"
Set-TextValue "D18" "2.431"
Set-TextValue "E18" "1.74%"
Set-TextValue "D19" "0.3403"
Set-TextValue "D20" "8.055"
Set-TextValue "E20" "-3.34%"
Set-TextValue "D21" "0.1370"
Set-TextValue "E21" "-2.10%"
Set-TextValue "E22" "0.79%"
Set-TextValue "D23" "0.001304"
Set-TextValue "E23" "3.15%"
Set-TextValue "D24" "0.004257"
Set-TextValue "E24" "9.27%"
Set-TextValue "E25" "3.97%"
Set-TextValue "E26" "-4.96%"
Set-TextValue "D38" "0.02724"
Set-TextValue "E38" "13.13%"
Set-TextValue "D39" "0.05763"
Set-TextValue "E39" "10.13%"
Set-TextValue "E40" "-6.70%"
Set-TextValue "D41" "0.007703"
Set-TextValue "E41" "-0.55%"
Set-TextValue "D42" "0.1445"
Set-TextValue "E42" "8.92%"
Set-TextValue "D43" "0.007696"
Set-TextValue "E43" "4.17%"
Set-TextValue "E44" "13.76%"
Set-TextValue "E45" "6.56%"
Set-TextValue "D46" "0.00006992"
Set-TextValue "E46" "12.11%"
Set-TextValue "E47" "0.02%"
Set-TextValue "D48" "0.05533"
Set-TextValue "E48" "20.95%"
Set-TextValue "E50" "0.02%"
Set-TextValue "E51" "0.02%"
